$d = $word.ActiveDocument

# 1) "quenching bath" -> "tempering bath" (appears twice)
$range = $d.Content
$range.Find.Execute("quenching bath", $false, $false, $false, $false, $false, $true, 1, $false, "tempering bath", 2)

# 2) Insert "take " into "Next, take a big cauldron ... tepid. &amp; three or four double handfuls"
$range2 = $d.Content
$range2.Find.Execute("three or four", $false, $false, $false, $false, $false, $true, 1, $false, "take three or four", 2)
